# Applies the "actividades" workbook edit:
#  - Hoja1 (sheet1): several activities were resolved/removed, one task was
#    renamed, one row's pairing was corrected, and the remaining rows were
#    re-packed (no gaps) -> dimension shrinks from A1:B14 to A1:B12.
#  - Modulo de Ventas (sheet2): four rows were appended at the bottom
#    (three activities that moved out of Hoja1 because they were finished,
#    plus one brand-new activity) -> dimension grows from A1:J26 to A1:J30.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Hoja1: clear out the old activity rows (2-14) and retype the final set.
# ---------------------------------------------------------------------
$ws1.Rows("2:14").Delete()

$sheet1Rows = @(
    @("Agregar Modulo Bitcoras", ""),
    @("Desarrollo Modulo  de Compras", ""),
    @("Desarrollar Modulo de Facturas", ""),
    @("Modulo de dashboard", "Desarrollo los querys del dashboard (1.-Cuanto llevan vendido las estaciones por diaen tiempo real, 2.-Top 5 de los productos vendidos  tiempo real , Top 5 de Clientes mensuales, Venta Total en tiempo real) en tiempo real me refiero que la query es al dia de hoy "),
    @("Modulo de Usuarios", "Cambiar al secuencia de los combos 1.- Alamcen , Sucursal  , Rol y hacerlos dependientes uno de cada uno es decir primero seleccionar Sucursal , Luego Almacen , y por ultimo el rol , no puede seleccionar ninguno otra opcion si no esta seleccionada una sucursal, tambien agregar un item que diga --Selecciona--"),
    @("Modulo de Usuarios", "Debe aparece el combo de sucursucal seleccionada la sucursal 1 o la predeterminda de solo lectura "),
    @("Modulo Agregar Estacion", ""),
    @("Modulo Descuentos", "Agregar la funcion de Agregar o descativar el Descuento "),
    @("Todo el sistema", "Los input type que son de numeros solo aceptar numeros con 2 decimales"),
    @("Todo el sistema", "Agregar loader en todas las peticiones por ajax"),
    @("Todos los Modulos", "agregar Async en todas las peticiones")
)

$r = 2
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "") {
        $ws1.Cells.Item($r, 2).Value = $row[1]
    }
    $r = $r + 1
}

$ws1.Range("A16").Select()

# ---------------------------------------------------------------------
# Modulo de Ventas: append the four new rows at the bottom of the list.
# ---------------------------------------------------------------------
$ws2.Range("A27").Value = "Calcular un rango en automatico del mayor hacia el infinito "
$ws2.Range("A28").Value = "Modulo Ventas: Agrega Campo para introducir cuanto le estan pagando y el sistema le diga cuando le corresponde de cambios"
$ws2.Range("A29").Value = "Modulo Ventas:Cuando agrego un item a la venta debe ser capaz de editar la cantidad (ya lo habiamos comentado este apartado)"
$ws2.Range("A30").Value = "Agregar PDF impresión de codigos de barras al modulo de productos"

$ws2.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 13
} catch {
}
$ws2.Range("A32").Select()
